# Apply rodada 27 Poisson-naive round/game-id updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple rows: update game id (A) and round number (E) ---
$ws.Range("A2").Value2 = 1
$ws.Range("E2").Value2 = 1
$ws.Range("A3").Value2 = 6
$ws.Range("E3").Value2 = 4
$ws.Range("A4").Value2 = 8
$ws.Range("E4").Value2 = 6
$ws.Range("A5").Value2 = 9
$ws.Range("E5").Value2 = 7
$ws.Range("A6").Value2 = 12
$ws.Range("E6").Value2 = 9
$ws.Range("A7").Value2 = 14
$ws.Range("E7").Value2 = 10
$ws.Range("A8").Value2 = 16
$ws.Range("E8").Value2 = 12
$ws.Range("A9").Value2 = 23
$ws.Range("E9").Value2 = 17
$ws.Range("A10").Value2 = 27
$ws.Range("E10").Value2 = 19
$ws.Range("A11").Value2 = 30
$ws.Range("E11").Value2 = 21
$ws.Range("A12").Value2 = 32
$ws.Range("E12").Value2 = 22
$ws.Range("A13").Value2 = 34
$ws.Range("E13").Value2 = 24
$ws.Range("A14").Value2 = 35
$ws.Range("E14").Value2 = 15
$ws.Range("A15").Value2 = 11
$ws.Range("E15").Value2 = 8
$ws.Range("A16").Value2 = 7
$ws.Range("E16").Value2 = 5
$ws.Range("A17").Value2 = 10
$ws.Range("E17").Value2 = 11
$ws.Range("A18").Value2 = 21
$ws.Range("E18").Value2 = 16
$ws.Range("A19").Value2 = 41
$ws.Range("E19").Value2 = 26
$ws.Range("A20").Value2 = 27
$ws.Range("E20").Value2 = 20
$ws.Range("A21").Value2 = 22
$ws.Range("E21").Value2 = 13
$ws.Range("A22").Value2 = 3
$ws.Range("E22").Value2 = 2
$ws.Range("A23").Value2 = 22
$ws.Range("E23").Value2 = 23
$ws.Range("A26").Value2 = 17
$ws.Range("E26").Value2 = 18
$ws.Range("A27").Value2 = 13
$ws.Range("E27").Value2 = 14

# --- Rows 24 & 25: full match data swapped between the two rows ---
# (row 24 now holds the match previously in row 25, and vice versa),
# with the game id (A) and round number (E) set for the new rodada 27 order.

# Row 24
$ws.Range("A24").Value2 = 24
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value2 = "2023-10-01"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value2 = "16:00"
$ws.Range("D24").Value2 = "Série A"
$ws.Range("E24").Value2 = 25
$ws.Range("F24").Value2 = "Sun"
$ws.Range("G24").Value2 = "Home"
$ws.Range("H24").Value2 = "D"
$ws.Range("I24").Value2 = 1
$ws.Range("J24").Value2 = 1
$ws.Range("K24").Value2 = "America (MG)"
$ws.Range("L24").Value2 = 0.6
$ws.Range("M24").Value2 = 1
$ws.Range("N24").Value2 = 59
$ws.Range("O24").Value2 = ""
$ws.Range("P24").Value2 = 17
$ws.Range("Q24").Value2 = 2
$ws.Range("R24").Value2 = 11.8
$ws.Range("S24").Value2 = 0.06
$ws.Range("T24").Value2 = 0.5
$ws.Range("U24").Value2 = 0
$ws.Range("V24").Value2 = 0
$ws.Range("W24").Value2 = 0
$ws.Range("X24").Value2 = 0.6
$ws.Range("Y24").Value2 = 0.03
$ws.Range("Z24").Value2 = 0.4
$ws.Range("AA24").Value2 = 0.4
$ws.Range("AB24").Value2 = 7
$ws.Range("AC24").Value2 = 6
$ws.Range("AD24").Value2 = 85.7
$ws.Range("AE24").Value2 = 0
$ws.Range("AF24").Value2 = 1.3
$ws.Range("AG24").Value2 = 0.3
$ws.Range("AH24").Value2 = 8317
$ws.Range("AI24").Value2 = 2814
$ws.Range("AJ24").Value2 = 1
$ws.Range("AK24").Value2 = 0.5
$ws.Range("AL24").Value2 = 0.8
$ws.Range("AM24").Value2 = 15
$ws.Range("AN24").Value2 = 39
$ws.Range("AO24").Value2 = 10
$ws.Range("AP24").Value2 = 4
$ws.Range("AQ24").Value2 = 51
$ws.Range("AR24").Value2 = 2
$ws.Range("AS24").Value2 = 1
$ws.Range("AT24").Value2 = 27
$ws.Range("AU24").Value2 = 10
$ws.Range("AV24").Value2 = 33
$ws.Range("AW24").Value2 = 2
$ws.Range("AX24").Value2 = 9
$ws.Range("AY24").Value2 = 4
$ws.Range("AZ24").Value2 = 10
$ws.Range("BA24").Value2 = 0
$ws.Range("BB24").Value2 = 12
$ws.Range("BC24").Value2 = 0
$ws.Range("BD24").Value2 = "Cruzeiro"

# Row 25
$ws.Range("A25").Value2 = 4
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value2 = "2023-04-29"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value2 = "18:30"
$ws.Range("D25").Value2 = "Série A"
$ws.Range("E25").Value2 = 3
$ws.Range("F25").Value2 = "Sat"
$ws.Range("G25").Value2 = "Home"
$ws.Range("H25").Value2 = "W"
$ws.Range("I25").Value2 = 3
$ws.Range("J25").Value2 = 2
$ws.Range("K25").Value2 = "America (MG)"
$ws.Range("L25").Value2 = 1.8
$ws.Range("M25").Value2 = 1.8
$ws.Range("N25").Value2 = 44
$ws.Range("O25").Value2 = 11844
$ws.Range("P25").Value2 = 11
$ws.Range("Q25").Value2 = 5
$ws.Range("R25").Value2 = 45.5
$ws.Range("S25").Value2 = 0.18
$ws.Range("T25").Value2 = 0.4
$ws.Range("U25").Value2 = 1
$ws.Range("V25").Value2 = 1
$ws.Range("W25").Value2 = 1
$ws.Range("X25").Value2 = 1
$ws.Range("Y25").Value2 = 0.1
$ws.Range("Z25").Value2 = 1.2
$ws.Range("AA25").Value2 = 1
$ws.Range("AB25").Value2 = 8
$ws.Range("AC25").Value2 = 6
$ws.Range("AD25").Value2 = 75
$ws.Range("AE25").Value2 = 0
$ws.Range("AF25").Value2 = 2.9
$ws.Range("AG25").Value2 = 0.9
$ws.Range("AH25").Value2 = 4196
$ws.Range("AI25").Value2 = 1827
$ws.Range("AJ25").Value2 = 2
$ws.Range("AK25").Value2 = 0.9
$ws.Range("AL25").Value2 = 0.7
$ws.Range("AM25").Value2 = 9
$ws.Range("AN25").Value2 = 24
$ws.Range("AO25").Value2 = 5
$ws.Range("AP25").Value2 = 1
$ws.Range("AQ25").Value2 = 25
$ws.Range("AR25").Value2 = 0
$ws.Range("AS25").Value2 = 1
$ws.Range("AT25").Value2 = 17
$ws.Range("AU25").Value2 = 8
$ws.Range("AV25").Value2 = 22
$ws.Range("AW25").Value2 = 5
$ws.Range("AX25").Value2 = 7
$ws.Range("AY25").Value2 = 3
$ws.Range("AZ25").Value2 = 11
$ws.Range("BA25").Value2 = 3
$ws.Range("BB25").Value2 = 10
$ws.Range("BC25").Value2 = 0
$ws.Range("BD25").Value2 = "Santos"
